# Update column F ("dSF") values to re-pulled data per commit message:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -7
$ws.Range("F5").Value = -3
$ws.Range("F8").Value = 4
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = -2
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = -2
